$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 32 and row 33 effectively swap their taxon data (with a couple of the
# B-column sort-order values changing to new numbers), per the diff.
# ---------------------------------------------------------------------------

# New row 32 values (previously row 33's taxon, plus updated B/Q/R values)
$ws.Cells.Item(32,1).Value2  = 112213305
$ws.Cells.Item(32,2).Value2  = 89503
$ws.Cells.Item(32,4).Value2  = "LC"
$ws.Cells.Item(32,5).Value2  = 5447
$ws.Cells.Item(32,6).Value2  = "Vedticka"
$ws.Cells.Item(32,7).Value2  = "Fuscoporia viticola"
$ws.Cells.Item(32,8).Value2  = "(Schwein.) Murrill"
$ws.Cells.Item(32,17).Value2 = 515748
$ws.Cells.Item(32,18).Value2 = 6704727

# New row 33 values (previously row 32's taxon, plus updated B/Q/R values)
$ws.Cells.Item(33,1).Value2  = 112213272
$ws.Cells.Item(33,2).Value2  = 89539
$ws.Cells.Item(33,4).Value2  = "NT"
$ws.Cells.Item(33,5).Value2  = 1202
$ws.Cells.Item(33,6).Value2  = "Ullticka"
$ws.Cells.Item(33,7).Value2  = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(33,8).Value2  = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(33,17).Value2 = 515738
$ws.Cells.Item(33,18).Value2 = 6704726

# ---------------------------------------------------------------------------
# New row 34
# ---------------------------------------------------------------------------

$ws.Cells.Item(34,1).Value2  = 112274505
$ws.Cells.Item(34,2).Value2  = 90800
$ws.Cells.Item(34,3).Value2  = "Ovaliderad"
$ws.Cells.Item(34,4).Value2  = "LC"
$ws.Cells.Item(34,5).Value2  = 4364
$ws.Cells.Item(34,6).Value2  = "Dropptaggsvamp"
$ws.Cells.Item(34,7).Value2  = "Hydnellum ferrugineum"
$ws.Cells.Item(34,8).Value2  = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(34,9).Value2  = "'2"
$ws.Cells.Item(34,10).Value2 = "fruktkroppar"
$ws.Cells.Item(34,16).Value2 = "Simsbodarna Ö, Dlr"
$ws.Cells.Item(34,17).Value2 = 515431
$ws.Cells.Item(34,18).Value2 = 6704883
$ws.Cells.Item(34,19).Value2 = 25
$ws.Cells.Item(34,20).Value2 = "Dalarna"
$ws.Cells.Item(34,21).Value2 = "Borlänge"
$ws.Cells.Item(34,22).Value2 = "Dalarna"
$ws.Cells.Item(34,23).Value2 = "Stora Tuna"
$ws.Cells.Item(34,25).Value2 = "'2023-09-23"
$ws.Cells.Item(34,26).Value2 = "'09:00"
$ws.Cells.Item(34,27).Value2 = "'2023-09-23"
$ws.Cells.Item(34,28).Value2 = "'10:40"
$ws.Cells.Item(34,30).Value2 = $false
$ws.Cells.Item(34,31).Value2 = $false
$ws.Cells.Item(34,33).Value2 = $false
$ws.Cells.Item(34,49).Value2 = "Håkan Sandin"
$ws.Cells.Item(34,50).Value2 = "Håkan Sandin"
